$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the existing AutoFilter criteria and unhide all the rows that were
#    hidden by the "RSSO" filter on column D.
$ws.ShowAllData()

# 2. Append the two new check rows to the bottom of the table.
$ws.Range("A143").Value = 247
$ws.Range("B143").Value = "SSH setup checks"
$ws.Range("C143").Value = "OS"
$ws.Range("D143").Value = "OS config"

$ws.Range("A144").Value = 248
# Leading apostrophe forces a text/quote-prefixed cell (matches style s="4"
# i.e. quotePrefix="1" used for this row in the target workbook).
$ws.Range("B144").Value = "'Pipeline library trailing space check"
$ws.Range("C144").Value = "Jenkins API"
$ws.Range("D144").Value = "Jenkins config"

# 3. Re-apply the AutoFilter over the (larger) A1:D152 range, with no filter
#    criteria selected this time, matching the committed workbook.
$ws.AutoFilterMode = $false
$ws.Range("A1:D152").AutoFilter()

# 4. The hidden _xlnm._FilterDatabase defined name should track the new
#    filter range.
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$D`$152"

# 5. Update the active selection / cursor position to match the author's
#    final position in the sheet.
$ws.Range("D145").Select()
